$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.537.73"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "2.246.19"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.14"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.61"
$ws.Range("E6").Value = "  -0.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  -1.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.68"
$ws.Range("E10").Value = "  -0.28%  "

$ws.Range("E11").Value = "  -1.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.19"
$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("E13").Value = "  +0.12%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.588.95"
$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.240.41"
$ws.Range("E15").Value = "  -4.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.831"
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.54"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").Value = "44.338.12"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").Value = "0.0₃0937"
$ws.Range("E19").Value = "  -2.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  -2.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.71"
$ws.Range("E21").Value = "  -2.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.19"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.04"
$ws.Range("E23").Value = "  -0.46%  "

$ws.Range("E24").Value = "  -5.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  -1.78%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.30"
$ws.Range("E27").Value = "  +4.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.74"
$ws.Range("E28").Value = "  -1.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.11"
$ws.Range("E29").Value = "  -3.17%  "

$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.92"
$ws.Range("E31").Value = "  -0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.04"
$ws.Range("E32").Value = "  -2.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0783"
$ws.Range("E33").Value = "  -1.08%  "

$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("E35").Value = "  -0.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.108"
$ws.Range("E36").Value = "  +1.82%  "

$ws.Range("E37").Value = "  -1.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.86"
$ws.Range("E38").Value = "  +6.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.23"
$ws.Range("E39").Value = "  +5.52%  "

$ws.Range("E40").Value = "  -5.54%  "

$ws.Range("E41").Value = "  -1.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0297"
$ws.Range("E42").Value = "  -0.71%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").Value = "1.803.76"
$ws.Range("E44").Value = "  +3.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.75"
$ws.Range("E45").Value = "  +10.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "81.96"
$ws.Range("E46").Value = "  -0.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.186"
$ws.Range("E47").Value = "  -2.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.33"
$ws.Range("E48").Value = "  -1.49%  "

$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.56"
$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.73"
$ws.Range("E51").Value = "  -1.53%  "
